# Apply cyclic rotation of species-observation data among rows 2, 3, 4, 5, 8.
# (row 6 and row 7 are untouched)
#
# The rotation moves data as follows (new row <- data that used to live in old row):
#   2 <- 4
#   3 <- 2
#   4 <- 8
#   5 <- 3
#   8 <- 5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the relevant (varying) column values for each source row, captured
# BEFORE any writes happen, so that writes to one row don't clobber data that
# still needs to be read for another row.
$cols = @("A","B","E","F","G","H","Q","R","AC")

$snapshot = @{}
foreach ($r in 2,3,4,5,8) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping: destination row -> source row (where the data used to be)
$rotation = @{
    2 = 4
    3 = 2
    4 = 8
    5 = 3
    8 = 5
}

foreach ($destRow in 2,3,4,5,8) {
    $srcRow = $rotation[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $val = $srcVals[$c]
        if ($c -eq "AC") {
            if ($null -eq $val -or $val -eq "") {
                $ws.Range("AC$destRow").Value = ""
            } else {
                $ws.Range("AC$destRow").Value = $val
            }
        } else {
            $ws.Range("$c$destRow").Value = $val
        }
    }
}
